# Added Page Object and WebDriverKeyowrds
# - Add a new "addEmployee" worksheet (after validLoginTest) with a small
#   Username/Password/Firstname/Lastname table of sample data.
# - Move the active selection off validLoginTest (to the new sheet), and
#   change validLoginTest's remembered selection to A1:B2.

$wb = $excel.ActiveWorkbook

# Update the selection remembered on "validLoginTest" before we move away
# from it (selecting a different range clears the previous activeCell).
$wsValid = $wb.Worksheets.Item("validLoginTest")
$wsValid.Range("A1:B2").Select() | Out-Null

# Add the new worksheet as the last tab in the workbook.
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "addEmployee"

# Header row
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "Firstname"
$ws.Range("D1").Value = "Lastname"

# Sample data row
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"
$ws.Range("C2").Value = "test fname"
$ws.Range("D2").Value = "test lname"

# Leave the new sheet's selection at F7, matching the authored workbook.
$ws.Range("F7").Select() | Out-Null
